# aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 227; Date = 44301; B = 2; C = 4; D = 115.5067860236789 },
    @{ Row = 228; Date = 44302; B = 0; C = 3; D = 86.63008951775916 },
    @{ Row = 229; Date = 44303; B = 0; C = 2; D = 57.75339301183945 }
)

# Column A uses a dedicated date style (border/font/alignment/number format);
# clone it from the last existing row rather than re-creating it from scratch.
$ws.Range("A226").Copy() | Out-Null

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
